$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy header formatting from F1 into new G1/H1, then set header text ---
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Cells.Item(1,6).Value = "Årsag"
$ws.Cells.Item(1,7).Value = "Ny leverandør"
$ws.Cells.Item(1,8).Value = "TCV_range"

# --- Data rows: set Reason (F), New supplier (G, only where applicable), TCV_range (H) ---
$ws.Cells.Item(2,6).Value = "Strategisk beslutning"
$ws.Cells.Item(2,8).Value = "40000-60000"
$ws.Cells.Item(3,6).Value = "Outsourcing af lønnen (anden leverandør)"
$ws.Cells.Item(3,7).Value = "Vil ikke oplyse"
$ws.Cells.Item(3,8).Value = "40000-60000"
$ws.Cells.Item(4,6).Value = "Strategisk beslutning"
$ws.Cells.Item(4,7).Value = "Lessor"
$ws.Cells.Item(4,8).Value = "40000-60000"
$ws.Cells.Item(5,6).Value = "Ikke oplyst"
$ws.Cells.Item(5,8).Value = "40000-60000"
$ws.Cells.Item(6,6).Value = "Fusionerer med anden virksomhed"
$ws.Cells.Item(6,8).Value = "40000-60000"
$ws.Cells.Item(7,6).Value = "Ikke oplyst"
$ws.Cells.Item(7,8).Value = "40000-60000"
$ws.Cells.Item(8,6).Value = "Ikke oplyst"
$ws.Cells.Item(8,8).Value = "40000-60000"
$ws.Cells.Item(9,6).Value = "Ikke oplyst"
$ws.Cells.Item(9,8).Value = "40000-60000"
$ws.Cells.Item(10,6).Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Cells.Item(10,7).Value = "Corp System"
$ws.Cells.Item(10,8).Value = "40000-60000"
$ws.Cells.Item(11,6).Value = "Ikke oplyst"
$ws.Cells.Item(11,8).Value = "40000-60000"
$ws.Cells.Item(12,6).Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Cells.Item(12,7).Value = "Zenegy"
$ws.Cells.Item(12,8).Value = "40000-60000"
$ws.Cells.Item(13,6).Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Cells.Item(13,8).Value = "40000-60000"
$ws.Cells.Item(14,6).Value = "Ikke oplyst"
$ws.Cells.Item(14,8).Value = "40000-60000"
$ws.Cells.Item(15,6).Value = "Ikke oplyst"
$ws.Cells.Item(15,8).Value = "40000-60000"
$ws.Cells.Item(16,6).Value = "Systemet (uddyb i bemærkninger)"
$ws.Cells.Item(16,8).Value = "40000-60000"
$ws.Cells.Item(17,6).Value = "Pris"
$ws.Cells.Item(17,7).Value = "Salary"
$ws.Cells.Item(17,8).Value = "40000-60000"
$ws.Cells.Item(18,6).Value = "Ikke oplyst"
$ws.Cells.Item(18,8).Value = "40000-60000"
$ws.Cells.Item(19,6).Value = "Ikke oplyst"
$ws.Cells.Item(19,8).Value = "40000-60000"
$ws.Cells.Item(20,6).Value = "Virksomheden lukker"
$ws.Cells.Item(20,8).Value = "40000-60000"
$ws.Cells.Item(21,6).Value = "Ikke flere medarbejdere i virksomheden"
$ws.Cells.Item(21,8).Value = "40000-60000"
$ws.Cells.Item(22,6).Value = "Ikke flere medarbejdere i virksomheden"
$ws.Cells.Item(22,8).Value = "40000-60000"
$ws.Cells.Item(23,6).Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Cells.Item(23,8).Value = "40000-60000"
$ws.Cells.Item(24,6).Value = "Pris"
$ws.Cells.Item(24,8).Value = "40000-60000"
$ws.Cells.Item(25,6).Value = "Ikke oplyst"
$ws.Cells.Item(25,8).Value = "40000-60000"
$ws.Cells.Item(26,6).Value = "Ikke oplyst"
$ws.Cells.Item(26,8).Value = "40000-60000"
$ws.Cells.Item(27,6).Value = "Ikke oplyst"
$ws.Cells.Item(27,8).Value = "40000-60000"
$ws.Cells.Item(28,6).Value = "Fusionerer med anden virksomhed"
$ws.Cells.Item(28,8).Value = "40000-60000"
$ws.Cells.Item(29,6).Value = "Systemet (uddyb i bemærkninger)"
$ws.Cells.Item(29,8).Value = "40000-60000"
$ws.Cells.Item(30,6).Value = "Virksomheden lukker"
$ws.Cells.Item(30,8).Value = "40000-60000"
$ws.Cells.Item(31,6).Value = "Pris"
$ws.Cells.Item(31,8).Value = "40000-60000"
$ws.Cells.Item(32,6).Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Cells.Item(32,8).Value = "40000-60000"
$ws.Cells.Item(33,6).Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Cells.Item(33,8).Value = "40000-60000"
$ws.Cells.Item(34,6).Value = "Ikke oplyst"
$ws.Cells.Item(34,7).Value = "DataLøn"
$ws.Cells.Item(34,8).Value = "40000-60000"
$ws.Cells.Item(35,6).Value = "Strategisk beslutning"
$ws.Cells.Item(35,8).Value = "40000-60000"
$ws.Cells.Item(36,6).Value = "Bruger ikke produktet"
$ws.Cells.Item(36,8).Value = "40000-60000"
